$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (ClaimMojio): mark Runmode as done
$ws.Range("C11").Value = "Y"

# Row 24 (ImportSIMs): mark Runmode as done
$ws.Range("C24").Value = "Y"

# Row 26 (ExportEvents): fill in Description/Runmode as done
$ws.Range("B26").Value = "Done"
$ws.Range("C26").Value = "N"

# New row 27: DeviceManager test case - copy formatting from row 26 first
$ws.Range("A26:C26").Copy($ws.Range("A27:C27"))
$ws.Range("A27").Value = "DeviceManager"
$ws.Range("B27").Value = ""
$ws.Range("C27").Value = "Y"

# Move active selection
$ws.Range("F24").Select()
